# Update scrape_timestamp values (column G) on the "raw" sheet from
# 2025-08-29T08:32:47Z to 2025-08-29T10:17:02Z for all data rows (2-241).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw")

$newTimestamp = "2025-08-29T10:17:02Z"

$lastRow = 241
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = $newTimestamp
}
